$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New trade row 4: copy styles from row 3 so date/boolean formatting (style index 1)
# carries over to the new row's A and G cells, then set the values.
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("G3").Copy()
$ws.Range("G4").PasteSpecial(-4122)

$ws.Range("A4").Value = 42641.539398148147
$ws.Range("B4").Value = $false
$ws.Range("C4").Value = 9835.02
$ws.Range("D4").Value = 9876.5
$ws.Range("E4").Value = 106.51
$ws.Range("F4").Value = 107.4
$ws.Range("G4").Value = $true
$ws.Range("H4").Value = 0.84
$ws.Range("I4").Value = $false
